{"js": "// Replace the 25 division-problem text values in the single 20x5 table.\n// Mapping is purely positional (document order), since some old values\n// repeat (e.g. \"64\u00f73=\" appears twice, and \"37\u00f72=\" is both a source and a\n// target value at different positions), so a global find/replace would be\n// ambiguous. We walk rows/cells in order instead.\nconst oldToNew = [\n  \"47\u00f77=\", \"89\u00f75=\",\n  \"66\u00f77=\", \"37\u00f72=\",\n  \"21\u00f77=\", \"28\u00f78=\",\n  \"68\u00f75=\", \"51\u00f75=\",\n  \"40\u00f73=\", \"40\u00f75=\",\n  \"67\u00f78=\", \"41\u00f72=\",\n  \"31\u00f74=\", \"79\u00f74=\",\n  \"33\u00f76=\", \"57\u00f74=\",\n  \"68\u00f79=\", \"74\u00f79=\",\n  \"64\u00f73=\", \"31\u00f78=\",\n  \"42\u00f75=\", \"44\u00f77=\",\n  \"64\u00f73=\", \"42\u00f74=\",\n  \"34\u00f79=\", \"64\u00f75=\",\n  \"81\u00f74=\", \"86\u00f72=\",\n  \"92\u00f79=\", \"37\u00f78=\",\n  \"37\u00f72=\", \"80\u00f74=\",\n  \"51\u00f77=\", \"44\u00f79=\",\n  \"20\u00f76=\", \"60\u00f77=\",\n  \"85\u00f72=\", \"93\u00f74=\",\n  \"73\u00f79=\", \"26\u00f73=\",\n  \"29\u00f77=\", \"20\u00f74=\",\n  \"27\u00f79=\", \"26\u00f76=\",\n  \"49\u00f72=\", \"52\u00f78=\",\n  \"24\u00f72=\", \"12\u00f79=\",\n  \"91\u00f74=\", \"68\u00f74=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Collect every non-empty paragraph (cell) in document order, along with\n// its current text, so we can match it against the expected old value\n// before replacing (defensive check against structural drift).\nconst targets = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    const cellBody = cell.body;\n    const paragraphs = cellBody.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    for (const paragraph of paragraphs.items) {\n      paragraph.load(\"text\");\n      targets.push(paragraph);\n    }\n  }\n}\nawait context.sync();\n\nlet pairIndex = 0;\nfor (const paragraph of targets) {\n  const text = paragraph.text;\n  if (text === \"\") {\n    continue; // blank spacer paragraphs/rows\n  }\n  if (pairIndex >= oldToNew.length) {\n    break;\n  }\n  const expectedOld = oldToNew[pairIndex];\n  const newValue = oldToNew[pairIndex + 1];\n  pairIndex += 2;\n\n  if (text !== expectedOld) {\n    throw new Error(\n      `Unexpected cell text \"${text}\" at position ${pairIndex / 2 - 1}; expected \"${expectedOld}\"`\n    );\n  }\n\n  const range = paragraph.getRange();\n  range.insertText(newValue, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem text values in the single 20x5 table.\n# Mapping is purely positional (document order / row-major cell order),\n# since some old values repeat (e.g. \"64\u00f73=\" appears twice, and \"37\u00f72=\"\n# is both a source and a target value at different positions), so a\n# global find/replace would be ambiguous. We walk the table cells in\n# row-major order instead, skipping the blank spacer rows.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$oldToNew = @(\n  \"47\u00f77=\", \"89\u00f75=\",\n  \"66\u00f77=\", \"37\u00f72=\",\n  \"21\u00f77=\", \"28\u00f78=\",\n  \"68\u00f75=\", \"51\u00f75=\",\n  \"40\u00f73=\", \"40\u00f75=\",\n  \"67\u00f78=\", \"41\u00f72=\",\n  \"31\u00f74=\", \"79\u00f74=\",\n  \"33\u00f76=\", \"57\u00f74=\",\n  \"68\u00f79=\", \"74\u00f79=\",\n  \"64\u00f73=\", \"31\u00f78=\",\n  \"42\u00f75=\", \"44\u00f77=\",\n  \"64\u00f73=\", \"42\u00f74=\",\n  \"34\u00f79=\", \"64\u00f75=\",\n  \"81\u00f74=\", \"86\u00f72=\",\n  \"92\u00f79=\", \"37\u00f78=\",\n  \"37\u00f72=\", \"80\u00f74=\",\n  \"51\u00f77=\", \"44\u00f79=\",\n  \"20\u00f76=\", \"60\u00f77=\",\n  \"85\u00f72=\", \"93\u00f74=\",\n  \"73\u00f79=\", \"26\u00f73=\",\n  \"29\u00f77=\", \"20\u00f74=\",\n  \"27\u00f79=\", \"26\u00f76=\",\n  \"49\u00f72=\", \"52\u00f78=\",\n  \"24\u00f72=\", \"12\u00f79=\",\n  \"91\u00f74=\", \"68\u00f74=\"\n)\n\n$pairIndex = 0\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # A cell's Range.Text always ends with the paragraph mark (\\r) plus the\n    # cell mark (\\x07); strip both so we can compare/replace just the\n    # visible text and detect genuinely empty spacer cells.\n    $text = $cellRange.Text.TrimEnd([char]13, [char]7)\n\n    if ($text -eq \"\") {\n      continue\n    }\n\n    $expectedOld = $oldToNew[$pairIndex]\n    $newValue = $oldToNew[$pairIndex + 1]\n    $pairIndex += 2\n\n    if ($text -ne $expectedOld) {\n      throw \"Unexpected cell text '$text' at row $r col $c; expected '$expectedOld'\"\n    }\n\n    $cellRange.Text = $newValue\n  }\n}\n"}
